# "listas de requesitos atualizada"
# Adds the RNF02..RNF10 non-functional requirements rows to the
# "Requesitos" sheet, continuing the existing table below row 29.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requesitos")

# Row 30 (old, empty placeholder row with its own one-off border style)
# gets overwritten with real content; rows 31-40 are new. Re-use the
# formatting already used by the existing empty placeholder rows
# (21-28) so we stay on the same style indices Excel already has,
# rather than inventing new ones.
$ws.Range("A28:D28").Copy() | Out-Null
$ws.Range("A30:D40").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Mirror the order the author actually typed things in: row 30 filled
# in one go (code, name, description), then the Cod. Requisito column
# typed downward for every remaining row, then each row's Nome /
# Descricao filled in afterwards. This keeps the shared-string table
# order identical to the authored workbook.
$ws.Cells.Item(30, 1).Value = "RNF02"
$ws.Cells.Item(30, 2).Value = "O tempo de execução"
$ws.Cells.Item(30, 3).Value = "O sistema deverá verificar se o usuário esta mais de 4 minutos logado sem nenhuma ação, o sistema devera após devera emitir um alerta de inativade."

$ws.Cells.Item(31, 1).Value = "RNF03"
$ws.Cells.Item(32, 1).Value = "RNF04"
$ws.Cells.Item(33, 1).Value = "RNF05"
$ws.Cells.Item(34, 1).Value = "RNF06"
$ws.Cells.Item(35, 1).Value = "RNF07"
$ws.Cells.Item(36, 1).Value = "RNF08"
$ws.Cells.Item(37, 1).Value = "RNF09"
$ws.Cells.Item(38, 1).Value = "RNF10"

$ws.Cells.Item(31, 2).Value = "O sistema operacional do sistema"
$ws.Cells.Item(31, 3).Value = "O sistema deverá ser acessado por varios tipos de sistemas operacionais ao mesmo tempo."

$ws.Cells.Item(32, 2).Value = "O desenvolvimento do sistema"
$ws.Cells.Item(32, 3).Value = "O sistema deverá ser desenvolvido ele toda em linguagem de programação JAVA"

$ws.Cells.Item(33, 2).Value = "Banco de dados do sistema"
$ws.Cells.Item(33, 3).Value = "O sistema devera se comunicar com o Banco de Dados Postgree"

$ws.Cells.Item(34, 2).Value = "Intregação com outro sistema"
$ws.Cells.Item(34, 3).Value = "O sistema deverá ser feito que possibilite a intregração com outro sistema."

$ws.Cells.Item(35, 2).Value = "O sistema deverá ter alta disponibilidade"
$ws.Cells.Item(35, 3).Value = "O sistema deverá ter alta disponibilidade, por exemplo 99% do tempo."

$ws.Cells.Item(36, 2).Value = "Tempo de processamento"
$ws.Cells.Item(36, 3).Value = "O sistema deverá processo N requisições por um determinado tempo."

$ws.Cells.Item(37, 2).Value = "Confiabilidade"
$ws.Cells.Item(37, 3).Value = "O sistema não deverá apresentar aos usuários quaisquer dados de cunho privativo."

# Row heights (row 38-40 keep the sheet's default row height).
$ws.Rows.Item(30).RowHeight = 60
$ws.Rows.Item(31).RowHeight = 45
$ws.Rows.Item(32).RowHeight = 30
$ws.Rows.Item(33).RowHeight = 30
$ws.Rows.Item(34).RowHeight = 30
$ws.Rows.Item(35).RowHeight = 30
$ws.Rows.Item(36).RowHeight = 30
$ws.Rows.Item(37).RowHeight = 30

# Rows 39-40 stay completely blank (already formatted via PasteSpecial above).

# Restore the view: scrolled down to the new rows, selection on B38.
$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 30
$excel.ActiveWindow.ScrollColumn = 1

# Page setup info present in the target file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
